{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\n// 1) Split the sentence ending \"...be implemented into the lights.\"\n//    into \"...be implemented into the lights\" + a new sentence describing\n//    the medicine-cabinet restock sensor.\nconst searchResults = context.document.body.search(\"be implemented into the lights.\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  const found = searchResults.items[0];\n  // Drop the trailing period and append the new sentence in its place.\n  found.insertText(\n    \"be implemented into the lights, along with a sensor in the medicine cabinet so it can be restocked whenever the stock is low.\",\n    Word.InsertLocation.replace\n  );\n}\n\n// 2) Remove the leftover \"_GoBack\" bookmark.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# $word / $app / $doc (ActiveDocument) are pre-seeded by the host.\n\n$d = $word.ActiveDocument\n\n# 1) Extend the sentence that currently ends \"...be implemented into the\n#    lights.\" with a new clause about a medicine-cabinet restock sensor.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"be implemented into the lights.\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"be implemented into the lights, along with a sensor in the medicine cabinet so it can be restocked whenever the stock is low.\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n$find.MatchCase = $true\n$find.Execute([ref]$find.Text, [ref]$find.MatchCase, $null, $null, $null, $null, [ref]$find.Forward, [ref]$find.Wrap, $null, [ref]$find.Replacement.Text, 2)  # wdReplaceAll\n\n# 2) Remove the leftover \"_GoBack\" bookmark.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
